# Refresh the crypto price/volume table (cryptos.xlsx) with the latest
# scraped figures, per the "Updated cryptos list ... with GitHub Actions"
# commit. All cells in this sheet are plain text (prices are formatted
# strings like "59.773.66", not real numbers), so each write forces the
# cell to Text format first and restores the default style afterwards —
# otherwise numeric-looking strings (e.g. "537.41") would silently be
# stored as real numbers by Excel's usual auto-detection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '59.773.66'
Set-TextValue 'E2' '  +1.05%  '
Set-TextValue 'D3' '2.634.38'
Set-TextValue 'E3' '  +1.76%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '537.41'
Set-TextValue 'E5' '  +1.56%  '
Set-TextValue 'E6' '  +3.32%  '
Set-TextValue 'E7' '  -0.04%  '
Set-TextValue 'E8' '  +0.37%  '
Set-TextValue 'D9' '6.54'
Set-TextValue 'E9' '  +1.91%  '
Set-TextValue 'E10' '  +2.28%  '
Set-TextValue 'E11' '  +1.52%  '
Set-TextValue 'E12' '  -1.66%  '
Set-TextValue 'D13' '3.106.31'
Set-TextValue 'E13' '  +2.11%  '
Set-TextValue 'D14' '59.695.64'
Set-TextValue 'E14' '  +1.08%  '
Set-TextValue 'D15' '20.93'
Set-TextValue 'E15' '  +2.11%  '
Set-TextValue 'B16' 'WrappedEther'
Set-TextValue 'C16' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D16' '2.627.32'
Set-TextValue 'E16' '  +1.31%  '
Set-TextValue 'B17' 'ShibaInu'
Set-TextValue 'C17' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D17' '0.0000135'
Set-TextValue 'E17' '  +1.46%  '
Set-TextValue 'D18' '341.98'
Set-TextValue 'E18' '  -0.83%  '
Set-TextValue 'E19' '  +2.64%  '
Set-TextValue 'E20' '  +1.40%  '
Set-TextValue 'D21' '6.40'
Set-TextValue 'E21' '  -0.39%  '
Set-TextValue 'D22' '0.999'
Set-TextValue 'D23' '67.57'
Set-TextValue 'E23' '  +0.33%  '
Set-TextValue 'E24' '  +1.70%  '
Set-TextValue 'E25' '  -0.54%  '
Set-TextValue 'E26' '  +0.09%  '
Set-TextValue 'E27' '  +2.59%  '
Set-TextValue 'D28' '0.0₃0753'
Set-TextValue 'E28' '  +4.86%  '
Set-TextValue 'E29' '  -0.03%  '
Set-TextValue 'E30' '  +3.87%  '
Set-TextValue 'E31' '  -0.06%  '
Set-TextValue 'D32' '18.96'
Set-TextValue 'E32' '  +1.36%  '
Set-TextValue 'D33' '150.79'
Set-TextValue 'E33' '  +1.31%  '
Set-TextValue 'D34' '4.01'
Set-TextValue 'E34' '  +1.42%  '
Set-TextValue 'E35' '  +2.14%  '
Set-TextValue 'B36' 'Stacks'
Set-TextValue 'C36' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D36' '1.46'
Set-TextValue 'E36' '  -1.21%  '
Set-TextValue 'B37' 'Fetch.AI'
Set-TextValue 'C37' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D37' '0.838'
Set-TextValue 'E37' '  +1.70%  '
Set-TextValue 'E38' '  +1.47%  '
Set-TextValue 'D39' '288.76'
Set-TextValue 'E39' '  +7.59%  '
Set-TextValue 'E40' '  +1.66%  '
Set-TextValue 'E41' '  -0.04%  '
Set-TextValue 'E42' '  +0.85%  '
Set-TextValue 'E43' '  -0.53%  '
Set-TextValue 'E44' '  -0.21%  '
Set-TextValue 'D45' '0.0533'
Set-TextValue 'E45' '  +3.76%  '
Set-TextValue 'D46' '1.967.62'
Set-TextValue 'E46' '  +0.42%  '
Set-TextValue 'E47' '  +1.58%  '
Set-TextValue 'D48' '18.49'
Set-TextValue 'E48' '  +1.80%  '
Set-TextValue 'E49' '  +2.72%  '
Set-TextValue 'D50' '110.88'
Set-TextValue 'E50' '  -0.59%  '
Set-TextValue 'D51' '4.74'
Set-TextValue 'E51' '  -0.23%  '
